# Banking Software Project Proposal - content edits
#
# 1) Merge the "I" / "mproved" / " fast and quick transactions" runs into a
#    single run reading "Improved fast and quick transactions".
# 2) Remove the stray _GoBack bookmark (bookmarkStart/bookmarkEnd) left over
#    from a previous save, leaving the (already empty) paragraph intact.
# 3) Split "Importance of securing you code." into
#    "Importance of securing you" + "r " + "code." so the sentence now
#    reads "Importance of securing your code."

$d = $word.ActiveDocument

# --- 1) Merge "I" + "mproved" + " fast and quick transactions" -----------
$d.Content.Find.Execute(
    "Improved fast and quick transactions",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Improved fast and quick transactions", 2) | Out-Null

# --- 2) Drop the leftover _GoBack bookmark --------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 3) "Importance of securing you code." -> "...your code." ------------
# Split the single space between "you" and "code." into its own run that
# carries the inserted "r", so the run boundaries end up as
# "Importance of securing you" | "r " | "code." (matching how Word keeps
# an in-place edit as separate runs once the revision is accepted).
$rng = $d.Content
$found = $rng.Find.Execute(
    "you code.",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found) {
    $wasTracking = $d.TrackRevisions
    $d.TrackRevisions = $true

    $spaceStart = $rng.Start + 3
    $spaceRange = $d.Range($spaceStart, $spaceStart + 1)
    $spaceRange.Text = "r "

    $d.TrackRevisions = $false
    while ($d.Revisions.Count -gt 0) {
        $d.Revisions(1).Accept() | Out-Null
    }
    $d.TrackRevisions = $wasTracking
}
